# Remove the "Upcountry District" column from the invoice detail sheet.
# The "Upcountry District" header/data lives in column H; deleting the
# entire column shifts the remaining columns (Upcountry Distance,
# Upcountry Amount) left by one, matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1").EntireColumn.Delete()
